$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text value corrections (typos / inconsistent data additions) ---
$ws.Range("F7").Value = "Programing"
$ws.Range("H2").Value = "Excell"
$ws.Range("C2").Value = "John12"
$ws.Range("D3").Value = "Programmers"
$ws.Range("J6").Value = "Communication_101"

# --- Number format changes (same dates, different display format) ---
$ws.Range("K2").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$ws.Range("G9").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

# --- Column width adjustment on column G ---
$ws.Columns("G").ColumnWidth = 22.5

# --- Selection moved ---
$ws.Range("O10").Select()
